# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted as row 97 of the "Rabanito"
# subset (Vega Central Mapocho de Santiago). Excel's native row-insert
# semantics push the existing rows 97..195 down to 98..196, which is
# exactly what the target diff shows (every row from 97 onward is the
# prior row's data, and the final row 195 spills into a brand new row 196).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 97; rows 97..195 shift down to 98..196
# and the sheet's used-range/dimension grows to R196 automatically.
$ws.Rows(97).Insert()

# Populate the newly inserted row 97 with the new observation.
$ws.Range("A97").Value = 9
$ws.Range("B97").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 44539
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = 300000001
$ws.Range("G97").Value = "Rabanito"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 7900
$ws.Range("K97").Value = 2500
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = 2747
$ws.Range("N97").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O97").Value = "Provincia de Chacabuco"
$ws.Range("P97").Value = 27
$ws.Range("Q97").Value = 100
$ws.Range("R97").Value = "Hortaliza"
